$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (Trigges Edita produtos / Criar): task finished -> Status goes from "?" (open) to "!" (done)
$ws.Range("F12").Value = "!"

# Row 15 (Teste): description + local were leftover placeholders from a previous
# evaluation ("Todas as tabelas e funções" / "T3 - Avalicação 2014-16") -> now
# finalized to reference the trigger work of this assignment
$ws.Range("B15").Value = "Todas as tabelas e trigger "
$ws.Range("C15").Value = "T4 - Trigges"

# Row 16 (Fazer / envio do trabalho): local also updated to the current assignment
$ws.Range("C16").Value = "T4 - Trigges"

# Update the active selection left on the sheet when the author saved
$ws.Range("C11").Select()
